$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "61.458.00"
$ws.Cells.Item(2, 5).Value = "  +0.65%  "
$ws.Cells.Item(3, 4).Value = "2.933.40"
$ws.Cells.Item(3, 5).Value = "  +0.38%  "
$ws.Cells.Item(4, 5).Value = "  -0.04%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "595.34"
$ws.Cells.Item(5, 5).Value = "  +0.81%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "145.19"
$ws.Cells.Item(6, 5).Value = "  +0.07%  "
$ws.Cells.Item(7, 5).Value = "  -0.04%  "
$ws.Cells.Item(8, 5).Value = "  -0.68%  "
$ws.Cells.Item(9, 5).Value = "  +1.79%  "
$ws.Cells.Item(10, 5).Value = "  -1.42%  "
$ws.Cells.Item(11, 5).Value = "  -0.54%  "
$ws.Cells.Item(12, 5).Value = "  -0.57%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "33.68"
$ws.Cells.Item(13, 5).Value = "  +0.03%  "
$ws.Cells.Item(14, 5).Value = "  +0.64%  "
$ws.Cells.Item(15, 4).Value = "3.420.00"
$ws.Cells.Item(15, 5).Value = "  +0.40%  "
$ws.Cells.Item(16, 4).Value = "61.410.47"
$ws.Cells.Item(16, 5).Value = "  +0.60%  "
$ws.Cells.Item(17, 5).Value = "  +0.18%  "
$ws.Cells.Item(18, 4).Value = "2.935.80"
$ws.Cells.Item(18, 5).Value = "  +0.59%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "432.31"
$ws.Cells.Item(19, 5).Value = "  +0.30%  "
$ws.Cells.Item(20, 5).Value = "  +0.12%  "
$ws.Cells.Item(21, 5).Value = "  -0.67%  "
$ws.Cells.Item(22, 5).Value = "  +0.51%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "81.85"
$ws.Cells.Item(23, 5).Value = "  +0.96%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "10.85"
$ws.Cells.Item(24, 5).Value = "  -1.57%  "
$ws.Cells.Item(25, 5).Value = "  -1.11%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "11.80"
$ws.Cells.Item(26, 5).Value = "  -2.10%  "
$ws.Cells.Item(27, 5).Value = "  -0.05%  "
$ws.Cells.Item(28, 5).Value = "  -3.13%  "
$ws.Cells.Item(29, 5).Value = "  -0.29%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "6.93"
$ws.Cells.Item(30, 5).Value = "  -2.65%  "
$ws.Cells.Item(31, 2).Value = "EthereumClassic"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "26.70"
$ws.Cells.Item(31, 5).Value = "  +0.67%  "
$ws.Cells.Item(32, 2).Value = "Hedera"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.110"
$ws.Cells.Item(32, 5).Value = "  +1.68%  "
$ws.Cells.Item(33, 5).Value = "  -0.02%  "
$ws.Cells.Item(34, 5).Value = "  +2.26%  "
$ws.Cells.Item(35, 5).Value = "  +0.09%  "
$ws.Cells.Item(36, 5).Value = "  +0.19%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "3.02"
$ws.Cells.Item(37, 5).Value = "  -1.25%  "
$ws.Cells.Item(38, 2).Value = "Kaspa"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.124"
$ws.Cells.Item(38, 5).Value = "  -0.46%  "
$ws.Cells.Item(39, 2).Value = "Stacks"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "2.01"
$ws.Cells.Item(39, 5).Value = "  +0.32%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "8.61"
$ws.Cells.Item(40, 5).Value = "  -0.01%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "42.66"
$ws.Cells.Item(41, 5).Value = "  +8.37%  "
$ws.Cells.Item(42, 5).Value = "  -1.07%  "
$ws.Cells.Item(43, 5).Value = "  -0.04%  "
$ws.Cells.Item(44, 4).Value = "2.706.93"
$ws.Cells.Item(44, 5).Value = "  -0.13%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "134.25"
$ws.Cells.Item(45, 5).Value = "  +1.99%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "365.05"
$ws.Cells.Item(46, 5).Value = "  -2.99%  "
$ws.Cells.Item(47, 5).Value = "  +0.02%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "23.79"
$ws.Cells.Item(48, 5).Value = "  -1.72%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "2.01"
$ws.Cells.Item(50, 5).Value = "  -1.43%  "
$ws.Cells.Item(51, 5).Value = "  -0.46%  "